$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename the header row labels: "_old" suffix -> "_FV2310", "_new" -> "_FV2404"
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Segmentname_FV2310"
$ws.Range("B1").Value = "Segmentgruppe_FV2310"
$ws.Range("C1").Value = "Segment_FV2310"
$ws.Range("D1").Value = "Datenelement_FV2310"
$ws.Range("E1").Value = "Segment ID_FV2310"
$ws.Range("F1").Value = "Code_FV2310"
$ws.Range("G1").Value = "Qualifier_FV2310"
$ws.Range("H1").Value = "Beschreibung_FV2310"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2310"
$ws.Range("J1").Value = "Bedingung_FV2310"
$ws.Range("K1").Value = "diff"
$ws.Range("L1").Value = "Segmentname_FV2404"
$ws.Range("M1").Value = "Segmentgruppe_FV2404"
$ws.Range("N1").Value = "Segment_FV2404"
$ws.Range("O1").Value = "Datenelement_FV2404"
$ws.Range("P1").Value = "Segment ID_FV2404"
$ws.Range("Q1").Value = "Code_FV2404"
$ws.Range("R1").Value = "Qualifier_FV2404"
$ws.Range("S1").Value = "Beschreibung_FV2404"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2404"
$ws.Range("U1").Value = "Bedingung_FV2404"

# ---------------------------------------------------------------------------
# 2) Convert A1:U65 into a native Excel Table ("Table1"), keeping the existing
#    header-row formatting intact (stash it, clear it so the table creation
#    does not capture a header dxf override, then restore it afterwards).
# ---------------------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("W1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("A1:U1").ClearFormats()

$rng = $ws.Range("A1:U65")
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

$ws.Range("W1").Copy()
$ws.Range("A1:U1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("W1").Clear()

# ---------------------------------------------------------------------------
# 3) Freeze the header row (split/freeze at row 2).
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
